# Auto-generated Excel COM-interop edit script
# Applies numeric 'want-to-go count' (F column) bumps across sheets 1-4,
# and reconciles the row 35-40 shuffle on sheet 4 (全部类型) caused by a newly
# inserted event (蔚蓝档案only) pushing subsequent rows down by one position.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# --- Sheet 1 (展览): update '想去人数' (F column) counts ---
$ws1.Range("F3").Value = 112
$ws1.Range("F7").Value = 1932
$ws1.Range("F8").Value = 5402
$ws1.Range("F9").Value = 1519
$ws1.Range("F10").Value = 154
$ws1.Range("F11").Value = 3100
$ws1.Range("F15").Value = 4280
$ws1.Range("F16").Value = 1025
$ws1.Range("F17").Value = 884
$ws1.Range("F18").Value = 1667
$ws1.Range("F19").Value = 2608
$ws1.Range("F21").Value = 25
$ws1.Range("F22").Value = 136
$ws1.Range("F24").Value = 975
$ws1.Range("F27").Value = 79
$ws1.Range("F29").Value = 1087
$ws1.Range("F30").Value = 383
$ws1.Range("F31").Value = 47
$ws1.Range("F32").Value = 147
$ws1.Range("F34").Value = 281
$ws1.Range("F36").Value = 1663
$ws1.Range("F37").Value = 2187
$ws1.Range("F38").Value = 1016
$ws1.Range("F41").Value = 608
$ws1.Range("F42").Value = 290
$ws1.Range("F44").Value = 651
$ws1.Range("F47").Value = 332
$ws1.Range("F49").Value = 135

# --- Sheet 2 (演出): update '想去人数' (F column) counts ---
$ws2.Range("F6").Value = 12
$ws2.Range("F10").Value = 145

# --- Sheet 3 (本地生活): update '想去人数' (F column) counts ---
$ws3.Range("F2").Value = 741

# --- Sheet 4 (全部类型): update '想去人数' (F column) counts (rows unaffected by the shuffle) ---
$ws4.Range("F2").Value = 741
$ws4.Range("F6").Value = 1932
$ws4.Range("F7").Value = 5402
$ws4.Range("F8").Value = 1519
$ws4.Range("F9").Value = 154
$ws4.Range("F11").Value = 3100
$ws4.Range("F14").Value = 4280
$ws4.Range("F15").Value = 1025
$ws4.Range("F16").Value = 1667
$ws4.Range("F17").Value = 12
$ws4.Range("F18").Value = 2608
$ws4.Range("F23").Value = 25
$ws4.Range("F25").Value = 145
$ws4.Range("F26").Value = 975
$ws4.Range("F29").Value = 79
$ws4.Range("F32").Value = 1087
$ws4.Range("F33").Value = 383
$ws4.Range("F34").Value = 47
$ws4.Range("F42").Value = 608
$ws4.Range("F43").Value = 290
$ws4.Range("F44").Value = 651
$ws4.Range("F46").Value = 332
$ws4.Range("F48").Value = 135

# --- Sheet 4 (全部类型): rows 35-40 shuffle ---
# A new event ('杭州·蔚蓝档案only') was inserted at row 35, shifting the former
# rows 35-39 down to 36-40, and the former row 40 event (青城山下) drops off the list.
# Each row's B:I cells are rewritten in place with their final target content.

$ws4.Range("B35").Value = '2024-06-01'
$ws4.Range("C35").Value = '杭州·蔚蓝档案only'
$ws4.Range("D35").Value = '北干街道萧杭路689号 杭州时尚外滩艺术中心'
$ws4.Range("E35").Value = '2024.06.01 09:00-06.01 18:00'
$ws4.Range("F35").Value = 147
$ws4.Range("G35").Value = 80
$ws4.Range("H35").Value = 'https://show.bilibili.com/platform/detail.html?id=84478'
$ws4.Range("I35").Value = '//i1.hdslb.com/bfs/openplatform/202404/z5lgl4tb1712719299126.jpeg'

$ws4.Range("B36").Value = '2024-06-08'
$ws4.Range("C36").Value = '杭州·第八届YH樱花动漫游戏文化节'
$ws4.Range("D36").Value = '德胜东路2539号 梦马汽车小镇'
$ws4.Range("E36").Value = '2024.06.08 10:00-06.10 17:00'
$ws4.Range("F36").Value = 1663
$ws4.Range("G36").Value = 65
$ws4.Range("H36").Value = 'https://show.bilibili.com/platform/detail.html?id=82687'
$ws4.Range("I36").Value = '//i0.hdslb.com/bfs/openplatform/202404/43sjLXZh1712910203022.jpeg'

$ws4.Range("B37").Value = '2024-06-09'
$ws4.Range("C37").Value = '杭州·第三届日夜国乙only'
$ws4.Range("D37").Value = '创意路1号 中国智谷富春园区'
$ws4.Range("E37").Value = '2024.06.09 10:00-06.09 23:00'
$ws4.Range("F37").Value = 2187
$ws4.Range("G37").Value = 58
$ws4.Range("H37").Value = 'https://show.bilibili.com/platform/detail.html?id=82618'
$ws4.Range("I37").Value = '//i2.hdslb.com/bfs/openplatform/202403/fXRzYEFH1710124366279.png'

$ws4.Range("B38").Value = '2024-06-15'
$ws4.Range("C38").Value = '杭州·次元盛典1.0'
$ws4.Range("D38").Value = '康候圣街99号 顺丰创新中心'
$ws4.Range("E38").Value = '2024.06.15 10:00-06.16 17:00'
$ws4.Range("F38").Value = 1016
$ws4.Range("G38").Value = 68
$ws4.Range("H38").Value = 'https://show.bilibili.com/platform/detail.html?id=83672'
$ws4.Range("I38").Value = '//i0.hdslb.com/bfs/openplatform/202404/yZAi07mM1712033477653.jpeg'

$ws4.Range("B39").Value = '2024-06-15'
$ws4.Range("C39").Value = '杭州·第三届动漫迷城嘉年华·毕业泳池'
$ws4.Range("D39").Value = '东新路21号 九龙仓君玺'
$ws4.Range("E39").Value = '2024.06.15 10:00-06.15 17:00'
$ws4.Range("F39").Value = 37
$ws4.Range("G39").Value = 70
$ws4.Range("H39").Value = 'https://show.bilibili.com/platform/detail.html?id=84338'
$ws4.Range("I39").Value = '//i1.hdslb.com/bfs/openplatform/202404/wQAlXTnK1713202337669.jpeg'

$ws4.Range("B40").Value = '2024-06-23'
$ws4.Range("C40").Value = '杭州·《亚米·跨越二次元》ACG经典动漫视听音乐会'
$ws4.Range("D40").Value = '金沙大道681号 金沙湖大剧院'
$ws4.Range("E40").Value = '2024.06.23 19:30-06.23 21:10'
$ws4.Range("F40").Value = 13
$ws4.Range("G40").Value = 80
$ws4.Range("H40").Value = 'https://show.bilibili.com/platform/detail.html?id=84041'
$ws4.Range("I40").Value = '//i2.hdslb.com/bfs/openplatform/202404/UhUuHfad1712564787267.jpeg'

